# Add "Wins", "Losses", "Ties" season-record columns (AD, AE, AF) to the
# player/roster sheet. Header row (row 1) gets the labels with the same
# bold/centered/bordered style used by the other header cells; every data
# row (2-53) gets the team's season record: 75 wins, 86 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, centered, bordered - style
# used by A1:AC1) onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2 through 53).
$wins = 75
$losses = 86
$ties = 0

for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
